# ventas.xlsx - "Add files via upload": append 20 new rows of comanda
# detail (rows 28-47) and fix the "Mesa" numbers on rows 26/27 so they
# are stored as real numbers instead of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix D26/D27: convert "30" stored as text to a real number 30 ---
$ws.Cells.Item(26, 4).Value = 30
$ws.Cells.Item(27, 4).Value = 30

# --- Append new comanda rows 28-47 ---
# Row 28
$ws.Cells.Item(28, 1).Value = 1705046158
$ws.Cells.Item(28, 2).Value = "2024-01-11 23:55:58"
$ws.Cells.Item(28, 5).Value = "Pizza"
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 7).Value = 500
$ws.Cells.Item(28, 8).Value = 500

# Row 29
$ws.Cells.Item(29, 1).Value = 1705046158
$ws.Cells.Item(29, 2).Value = "2024-01-11 23:55:58"
$ws.Cells.Item(29, 5).Value = "Hamburguesas"
$ws.Cells.Item(29, 6).Value = 1
$ws.Cells.Item(29, 7).Value = 300
$ws.Cells.Item(29, 8).Value = 300

# Row 30
$ws.Cells.Item(30, 1).Value = 1705046158
$ws.Cells.Item(30, 2).Value = "2024-01-11 23:55:58"
$ws.Cells.Item(30, 5).Value = "Crostinis"
$ws.Cells.Item(30, 6).Value = 1
$ws.Cells.Item(30, 7).Value = 50
$ws.Cells.Item(30, 8).Value = 50

# Row 31
$ws.Cells.Item(31, 1).Value = 1705046158
$ws.Cells.Item(31, 2).Value = "2024-01-11 23:55:58"
$ws.Cells.Item(31, 5).Value = "Cerveza"
$ws.Cells.Item(31, 6).Value = 1
$ws.Cells.Item(31, 7).Value = 50
$ws.Cells.Item(31, 8).Value = 50

# Row 32
$ws.Cells.Item(32, 1).Value = 1705046215
$ws.Cells.Item(32, 2).Value = "2024-01-11 23:56:55"
$ws.Cells.Item(32, 3).Value = "Victor Hugo"
$ws.Cells.Item(32, 4).Value = 15
$ws.Cells.Item(32, 5).Value = "Pizza"
$ws.Cells.Item(32, 6).Value = 1
$ws.Cells.Item(32, 7).Value = 500
$ws.Cells.Item(32, 8).Value = 500

# Row 33
$ws.Cells.Item(33, 1).Value = 1705046215
$ws.Cells.Item(33, 2).Value = "2024-01-11 23:56:55"
$ws.Cells.Item(33, 3).Value = "Victor Hugo"
$ws.Cells.Item(33, 4).Value = 15
$ws.Cells.Item(33, 5).Value = "Hamburguesas"
$ws.Cells.Item(33, 6).Value = 1
$ws.Cells.Item(33, 7).Value = 300
$ws.Cells.Item(33, 8).Value = 300

# Row 34
$ws.Cells.Item(34, 1).Value = 1705046215
$ws.Cells.Item(34, 2).Value = "2024-01-11 23:56:55"
$ws.Cells.Item(34, 3).Value = "Victor Hugo"
$ws.Cells.Item(34, 4).Value = 15
$ws.Cells.Item(34, 5).Value = "Crostinis"
$ws.Cells.Item(34, 6).Value = 1
$ws.Cells.Item(34, 7).Value = 50
$ws.Cells.Item(34, 8).Value = 50

# Row 35
$ws.Cells.Item(35, 1).Value = 1705046215
$ws.Cells.Item(35, 2).Value = "2024-01-11 23:56:55"
$ws.Cells.Item(35, 3).Value = "Victor Hugo"
$ws.Cells.Item(35, 4).Value = 15
$ws.Cells.Item(35, 5).Value = "Cerveza"
$ws.Cells.Item(35, 6).Value = 1
$ws.Cells.Item(35, 7).Value = 50
$ws.Cells.Item(35, 8).Value = 50

# Row 36
$ws.Cells.Item(36, 1).Value = 1705048142
$ws.Cells.Item(36, 2).Value = "2024-01-12 00:29:02"
$ws.Cells.Item(36, 5).Value = "Pizza"
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(36, 7).Value = 500
$ws.Cells.Item(36, 8).Value = 500

# Row 37
$ws.Cells.Item(37, 1).Value = 1705048142
$ws.Cells.Item(37, 2).Value = "2024-01-12 00:29:02"
$ws.Cells.Item(37, 5).Value = "Hamburguesas"
$ws.Cells.Item(37, 6).Value = 1
$ws.Cells.Item(37, 7).Value = 300
$ws.Cells.Item(37, 8).Value = 300

# Row 38
$ws.Cells.Item(38, 1).Value = 1705048142
$ws.Cells.Item(38, 2).Value = "2024-01-12 00:29:02"
$ws.Cells.Item(38, 5).Value = "Crostinis"
$ws.Cells.Item(38, 6).Value = 1
$ws.Cells.Item(38, 7).Value = 50
$ws.Cells.Item(38, 8).Value = 50

# Row 39
$ws.Cells.Item(39, 1).Value = 1705048142
$ws.Cells.Item(39, 2).Value = "2024-01-12 00:29:02"
$ws.Cells.Item(39, 5).Value = "Cerveza"
$ws.Cells.Item(39, 6).Value = 1
$ws.Cells.Item(39, 7).Value = 50
$ws.Cells.Item(39, 8).Value = 50

# Row 40
$ws.Cells.Item(40, 1).Value = 1705048205
$ws.Cells.Item(40, 2).Value = "2024-01-12 00:30:05"
$ws.Cells.Item(40, 5).Value = "Pizza"
$ws.Cells.Item(40, 6).Value = 1
$ws.Cells.Item(40, 7).Value = 500
$ws.Cells.Item(40, 8).Value = 500

# Row 41
$ws.Cells.Item(41, 1).Value = 1705048205
$ws.Cells.Item(41, 2).Value = "2024-01-12 00:30:05"
$ws.Cells.Item(41, 5).Value = "Hamburguesas"
$ws.Cells.Item(41, 6).Value = 1
$ws.Cells.Item(41, 7).Value = 300
$ws.Cells.Item(41, 8).Value = 300

# Row 42
$ws.Cells.Item(42, 1).Value = 1705048205
$ws.Cells.Item(42, 2).Value = "2024-01-12 00:30:05"
$ws.Cells.Item(42, 5).Value = "Crostinis"
$ws.Cells.Item(42, 6).Value = 1
$ws.Cells.Item(42, 7).Value = 50
$ws.Cells.Item(42, 8).Value = 50

# Row 43
$ws.Cells.Item(43, 1).Value = 1705048205
$ws.Cells.Item(43, 2).Value = "2024-01-12 00:30:05"
$ws.Cells.Item(43, 5).Value = "Cerveza"
$ws.Cells.Item(43, 6).Value = 1
$ws.Cells.Item(43, 7).Value = 50
$ws.Cells.Item(43, 8).Value = 50

# Row 44
$ws.Cells.Item(44, 1).Value = 1705048668
$ws.Cells.Item(44, 2).Value = "2024-01-12 00:37:48"
$ws.Cells.Item(44, 5).Value = "Pizza"
$ws.Cells.Item(44, 6).Value = 1
$ws.Cells.Item(44, 7).Value = 500
$ws.Cells.Item(44, 8).Value = 500

# Row 45
$ws.Cells.Item(45, 1).Value = 1705048668
$ws.Cells.Item(45, 2).Value = "2024-01-12 00:37:48"
$ws.Cells.Item(45, 5).Value = "Hamburguesas"
$ws.Cells.Item(45, 6).Value = 4
$ws.Cells.Item(45, 7).Value = 300
$ws.Cells.Item(45, 8).Value = 1200

# Row 46
$ws.Cells.Item(46, 1).Value = 1705048668
$ws.Cells.Item(46, 2).Value = "2024-01-12 00:37:48"
$ws.Cells.Item(46, 5).Value = "Crostinis"
$ws.Cells.Item(46, 6).Value = 3
$ws.Cells.Item(46, 7).Value = 50
$ws.Cells.Item(46, 8).Value = 150

# Row 47
$ws.Cells.Item(47, 1).Value = 1705048668
$ws.Cells.Item(47, 2).Value = "2024-01-12 00:37:48"
$ws.Cells.Item(47, 5).Value = "Cerveza"
$ws.Cells.Item(47, 6).Value = 4
$ws.Cells.Item(47, 7).Value = 50
$ws.Cells.Item(47, 8).Value = 200

